# Updates cryptos list values (price & 1h volume change) scraped on
# Mon Feb  5 14:45:36 UTC 2024, including a rank swap between
# EnergySwap and FraxShare (rows 43/44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.321.58"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "2.324.24"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'302.29"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "'98.02"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.507"
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.503"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "'35.56"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "'19.56"
$ws.Range("E11").Value = "  +7.59%  "
$ws.Range("D12").Value = "'0.0796"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "'6.90"
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "2.689.14"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").Value = "2.328.77"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").Value = "'0.792"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "43.246.83"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "'12.71"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "0.0₃0899"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "'6.07"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "'68.01"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'237.28"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("E24").Value = "  +4.41%  "
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'25.06"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("D29").Value = "'164.28"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").Value = "'9.13"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").Value = "'33.13"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'5.01"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").Value = "'17.88"
$ws.Range("E34").Value = "  +3.18%  "
$ws.Range("D35").Value = "'4.51"
$ws.Range("E35").Value = "  -6.41%  "
$ws.Range("D36").Value = "'0.0705"
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("D40").Value = "'2.79"
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("D41").Value = "'0.109"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "1.984.62"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'10.66"
$ws.Range("E43").Value = "  +5.97%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'18.98"
$ws.Range("E44").Value = "  +5.93%  "
$ws.Range("D46").Value = "'2.07"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "'2.80"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "2.555.09"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").Value = "'2.86"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").Value = "'53.92"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").Value = "'72.66"
$ws.Range("E51").Value = "  +0.71%  "

# Excel marks cells entered with a leading apostrophe (used above to force
# text storage for values that would otherwise be auto-parsed as numbers)
# with a quotePrefix style. Restore the default "Normal" style on those
# cells so their formatting matches the rest of the untouched column.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
